$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update solver/model parameters (B6:B9)
$ws.Range("B6").Value = 0.58595469351715701
$ws.Range("B7").Value = 176.4510567799318
$ws.Range("B8").Value = 1.0393701324990459
$ws.Range("B9").Value = 0.67944147424124168

# Update observed confirmed-case counts (F column), shifted by one day with new last value
$ws.Range("F15").Value = 72
$ws.Range("F16").Value = 117
$ws.Range("F17").Value = 130
$ws.Range("F18").Value = 188
$ws.Range("F19").Value = 240
$ws.Range("F20").Value = 351
$ws.Range("F21").Value = 670
$ws.Range("F22").Value = 824
$ws.Range("F23").ClearContents()

# Error term (K23) no longer has data to compare against, clear it
$ws.Range("K23").ClearContents()

# Expand the sum-of-squared-error range to include the new last data point
$ws.Range("D7").Formula = "=SUM(K15:K22)"

# Move the active selection like the author's last edit
$ws.Range("P7").Select()
